# "updated dummy data with actual data"
# Adds three new attribute columns (InputType, ValidationConstraint, MaxLength)
# to the User entity properties sheet header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "InputType"
$ws.Range("F1").Value = "ValidationConstraint"
$ws.Range("G1").Value = "MaxLength"

# ColumnWidth is stored/round-tripped in "character" units with a pixel-rounding
# step, so nudge the requested width down slightly to land on the clean target
# value (12 / 18 / 10) once persisted back to OOXML.
$ws.Columns.Item(5).ColumnWidth = 11.1
$ws.Columns.Item(6).ColumnWidth = 17.1
$ws.Columns.Item(7).ColumnWidth = 9.1

$ws.Range("G1").Select()
